# feat：update skill icon excel
# Populate the (previously empty) H/I "icon id" columns for skill rows
# that were missing selectIconGuid / changeIconGuid values, and move the
# active selection to I46 as recorded in the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (H = selectIconGuid, I = changeIconGuid)
$iconData = @(
    @{ Row = 8;  H = 337784; I = 337785 },
    @{ Row = 15; H = 326721; I = 326720 },
    @{ Row = 16; H = 326711; I = 326709 },
    @{ Row = 17; H = 326713; I = 326716 },
    @{ Row = 18; H = 337503; I = 337502 },
    @{ Row = 19; H = 326714; I = 326717 },
    @{ Row = 20; H = 295333; I = 295325 },
    @{ Row = 21; H = 326696; I = 326700 },
    @{ Row = 22; H = 326694; I = 326702 },
    @{ Row = 23; H = 326697; I = 326692 },
    @{ Row = 24; H = 326707; I = 326715 },
    @{ Row = 25; H = 327134; I = 327105 },
    @{ Row = 26; H = 295346; I = 295338 },
    @{ Row = 27; H = 327131; I = 327130 },
    @{ Row = 28; H = 327106; I = 327109 },
    @{ Row = 29; H = 327135; I = 327117 },
    @{ Row = 30; H = 327108; I = 327113 },
    @{ Row = 31; H = 327125; I = 327124 },
    @{ Row = 32; H = 327107; I = 327102 },
    @{ Row = 33; H = 327132; I = 327127 },
    @{ Row = 34; H = 327110; I = 327103 },
    @{ Row = 35; H = 327092; I = 327174 },
    @{ Row = 36; H = 327111; I = 327084 },
    @{ Row = 37; H = 295328; I = 295327 },
    @{ Row = 38; H = 327086; I = 327093 },
    @{ Row = 39; H = 327096; I = 327095 },
    @{ Row = 40; H = 327169; I = 327168 },
    @{ Row = 41; H = 327091; I = 327082 },
    @{ Row = 42; H = 327154; I = 327156 },
    @{ Row = 43; H = 327090; I = 327085 },
    @{ Row = 44; H = 327112; I = 327098 }
)

foreach ($item in $iconData) {
    $ws.Cells.Item($item.Row, 8).Value = $item.H   # column H
    $ws.Cells.Item($item.Row, 9).Value = $item.I   # column I
}

# Match the workbook's recorded view state after the edit: scrolled down
# and the active cell/selection moved to I46.
$ws.Range("I46").Select()
